# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C text swaps (coin list reordering) ---
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("B43").Value = 'Arweave'
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

# --- Column D/E numeric-looking text updates (Price / Volume(1h)) ---
# Values such as "0.999" or "2.70" would be auto-coerced to numbers by a plain
# Range.Value assignment (losing the intended text formatting / trailing zeros),
# so we stage each value as a text FORMULA result (="...") in a scratch cell far
# off-sheet, copy it, and PasteSpecial with xlPasteValues (-4163) into the
# destination - this keeps the cell as text (same as the original inlineStr)
# without altering any cell styles or leaving a formula behind.

$ws.Range("Z1").Formula = '="61.972.37"'
$ws.Range("AA1").Formula = '="  -0.93%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="3.001.96"'
$ws.Range("AA1").Formula = '="  -0.33%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.999"'
$ws.Range("AA1").Formula = '="  -0.02%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="593.28"'
$ws.Range("AA1").Formula = '="  +1.40%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="146.96"'
$ws.Range("AA1").Formula = '="  +0.28%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -0.03%  "'
$ws.Range("Z1").Copy()
$ws.Range("E7").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="2.999.42"'
$ws.Range("AA1").Formula = '="  -0.39%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -2.18%  "'
$ws.Range("Z1").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +6.73%  "'
$ws.Range("Z1").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -0.32%  "'
$ws.Range("Z1").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.456"'
$ws.Range("AA1").Formula = '="  -0.86%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +0.39%  "'
$ws.Range("Z1").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="34.28"'
$ws.Range("AA1").Formula = '="  -1.28%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +2.92%  "'
$ws.Range("Z1").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="3.495.11"'
$ws.Range("AA1").Formula = '="  -0.29%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="6.97"'
$ws.Range("AA1").Formula = '="  -1.82%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="61.910.30"'
$ws.Range("AA1").Formula = '="  -0.91%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="3.000.95"'
$ws.Range("AA1").Formula = '="  -0.33%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="446.90"'
$ws.Range("AA1").Formula = '="  -2.66%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +0.91%  "'
$ws.Range("Z1").Copy()
$ws.Range("E21").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.685"'
$ws.Range("AA1").Formula = '="  -0.66%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -1.03%  "'
$ws.Range("Z1").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="82.24"'
$ws.Range("AA1").Formula = '="  +0.68%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="10.98"'
$ws.Range("AA1").Formula = '="  +9.28%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="2.22"'
$ws.Range("AA1").Formula = '="  -0.08%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="12.12"'
$ws.Range("AA1").Formula = '="  -2.09%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +0.06%  "'
$ws.Range("Z1").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="2.70"'
$ws.Range("AA1").Formula = '="  +2.94%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +0.04%  "'
$ws.Range("Z1").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +2.05%  "'
$ws.Range("Z1").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -0.49%  "'
$ws.Range("Z1").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="27.44"'
$ws.Range("AA1").Formula = '="  -2.18%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.111"'
$ws.Range("AA1").Formula = '="  +1.26%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.0₃0845"'
$ws.Range("AA1").Formula = '="  +4.11%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -0.52%  "'
$ws.Range("Z1").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  +0.39%  "'
$ws.Range("Z1").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="50.13"'
$ws.Range("AA1").Formula = '="  -0.43%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="2.04"'
$ws.Range("AA1").Formula = '="  -3.93%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="8.97"'
$ws.Range("AA1").Formula = '="  -1.74%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="2.92"'
$ws.Range("AA1").Formula = '="  +0.52%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.123"'
$ws.Range("AA1").Formula = '="  +3.87%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="40.74"'
$ws.Range("AA1").Formula = '="  +8.58%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="395.38"'
$ws.Range("AA1").Formula = '="  +1.39%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.279"'
$ws.Range("AA1").Formula = '="  +3.52%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="0.0350"'
$ws.Range("AA1").Formula = '="  -2.43%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="2.714.43"'
$ws.Range("AA1").Formula = '="  -0.80%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="132.58"'
$ws.Range("AA1").Formula = '="  +2.39%  "'
$ws.Range("Z1:AA1").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -1.80%  "'
$ws.Range("Z1").Copy()
$ws.Range("E50").PasteSpecial(-4163)
$ws.Range("Z1").Formula = '="  -1.80%  "'
$ws.Range("Z1").Copy()
$ws.Range("E51").PasteSpecial(-4163)

# --- Clean up the scratch cells ---
$ws.Range("Z1:AA1").ClearContents()
$excel.CutCopyMode = 0
